$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.220.29'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.860.51'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''236.40'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '''0.4673'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '''0.2857'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('D9').Value = '''0.06525'
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').Value = '''21.87'
$ws.Range('E10').Value = '  +8.31%  '
$ws.Range('D11').Value = '''0.07913'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').Value = '''97.23'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '1.872.00'
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').Value = '''5.153'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').Value = '''0.6802'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '''269.55'
$ws.Range('E16').Value = '  -3.89%  '
$ws.Range('D17').Value = '30.228.31'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '''13.51'
$ws.Range('E18').Value = '  +6.80%  '
$ws.Range('D20').Value = '''0.000007337'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').Value = '2.113.88'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').Value = '''5.322'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('D23').Value = '''1.002'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '''6.167'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').Value = '''167.52'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = '''9.200'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').Value = '''18.87'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').Value = '''1.949'
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('D29').Value = '''1.384'
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('D30').Value = '''0.09841'
$ws.Range('E30').Value = '  +2.53%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.480'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''4.348'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('D33').Value = '''4.050'
$ws.Range('E33').Value = '  -1.69%  '
$ws.Range('D34').Value = '''0.04702'
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('D35').Value = '''1.131'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').Value = '''0.7001'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').Value = '''0.01870'
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').Value = '''2.625'
$ws.Range('E39').Value = '  +4.07%  '
$ws.Range('D40').Value = '''75.49'
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('D41').Value = '''6.266'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').Value = '''1.940'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = '''0.8516'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '''0.4154'
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').Value = '''103.07'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '''952.34'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''7.157'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').Value = '''9.233'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').Value = '''34.09'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').Value = '''0.05649'
$ws.Range('E51').Value = '  +0.33%  '
